$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was FAPs target; becomes ECs target with updated stats)
$ws.Cells.Item(2,1).Value  = "FAPs"
$ws.Cells.Item(2,2).Value  = "Wnt2"
$ws.Cells.Item(2,3).Value  = "Fzd2"
$ws.Cells.Item(2,4).Value  = "ECs"
$ws.Cells.Item(2,5).Value  = 3
$ws.Cells.Item(2,6).Value  = 1
$ws.Cells.Item(2,7).Value  = 0.8330250000000001
$ws.Cells.Item(2,8).Value  = 2.499075
$ws.Cells.Item(2,9).Value  = 1
$ws.Cells.Item(2,10).Value = 1
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.07629866666666667
$ws.Cells.Item(2,14).Value = 0.228896
$ws.Cells.Item(2,15).Value = 0.004108848954870246
$ws.Cells.Item(2,16).Value = 0.004108848954870246
$ws.Cells.Item(2,17).Value = 0.06355869680000001
$ws.Cells.Item(2,18).Value = 0.5720282712000001
$ws.Cells.Item(2,19).Value = 0.004108848954870246
$ws.Cells.Item(2,20).Value = 0.004108848954870246

# Row 3 (was sCs target; becomes FAPs target with updated stats)
$ws.Cells.Item(3,1).Value  = "FAPs"
$ws.Cells.Item(3,2).Value  = "Wnt2"
$ws.Cells.Item(3,3).Value  = "Fzd2"
$ws.Cells.Item(3,4).Value  = "FAPs"
$ws.Cells.Item(3,5).Value  = 3
$ws.Cells.Item(3,6).Value  = 1
$ws.Cells.Item(3,7).Value  = 0.8330250000000001
$ws.Cells.Item(3,8).Value  = 2.499075
$ws.Cells.Item(3,9).Value  = 1
$ws.Cells.Item(3,10).Value = 1
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 12.99468133333333
$ws.Cells.Item(3,14).Value = 38.984044
$ws.Cells.Item(3,15).Value = 0.6997918200668237
$ws.Cells.Item(3,16).Value = 0.6997918200668237
$ws.Cells.Item(3,17).Value = 10.8248944177
$ws.Cells.Item(3,18).Value = 97.4240497593
$ws.Cells.Item(3,19).Value = 0.6997918200668237
$ws.Cells.Item(3,20).Value = 0.6997918200668237

# Row 4 (new row; sCs target)
$ws.Cells.Item(4,1).Value  = "FAPs"
$ws.Cells.Item(4,2).Value  = "Wnt2"
$ws.Cells.Item(4,3).Value  = "Fzd2"
$ws.Cells.Item(4,4).Value  = "sCs"
$ws.Cells.Item(4,5).Value  = 3
$ws.Cells.Item(4,6).Value  = 1
$ws.Cells.Item(4,7).Value  = 0.8330250000000001
$ws.Cells.Item(4,8).Value  = 2.499075
$ws.Cells.Item(4,9).Value  = 1
$ws.Cells.Item(4,10).Value = 1
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.498373000000001
$ws.Cells.Item(4,14).Value = 16.495119
$ws.Cells.Item(4,15).Value = 0.2960993309783061
$ws.Cells.Item(4,16).Value = 0.2960993309783061
$ws.Cells.Item(4,17).Value = 4.580282168325001
$ws.Cells.Item(4,18).Value = 41.22253951492501
$ws.Cells.Item(4,19).Value = 0.2960993309783061
$ws.Cells.Item(4,20).Value = 0.2960993309783061
